$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the 4 new rows that appear in the updated data ---
# New row for WBA (becomes row 37) and Coventry (becomes row 38),
# both inserted before the existing "Southampton" row (currently row 37).
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).Insert()

# New row for Middlesbrough, inserted before the existing "Oxford" row
# (currently row 39, after the two inserts above it is row 41).
$ws.Rows.Item(41).Insert()

# New row for Sheff Utd, inserted before the existing "Swansea" row
# (currently row 40, after the three inserts above it is row 43).
$ws.Rows.Item(43).Insert()

# --- Copy the formatting (bold/border/center) used by column A onto the
#     A cells of the newly inserted rows, matching the rest of the sheet ---
$ws.Range("A2").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A41").PasteSpecial(-4122)
$ws.Range("A43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Populate the new rows with their data ---
$ws.Range("A37").Value2 = 35
$ws.Range("B37").Value2 = 1908387
$ws.Range("C37").Value2 = "Watford"
$ws.Range("D37").Value2 = "WBA"

$ws.Range("A38").Value2 = 36
$ws.Range("B38").Value2 = 1908400
$ws.Range("C38").Value2 = "Watford"
$ws.Range("D38").Value2 = "Coventry"

$ws.Range("A41").Value2 = 39
$ws.Range("B41").Value2 = 1908455
$ws.Range("C41").Value2 = "Watford"
$ws.Range("D41").Value2 = "Middlesbrough"

$ws.Range("A43").Value2 = 41
$ws.Range("B43").Value2 = 1908486
$ws.Range("C43").Value2 = "Watford"
$ws.Range("D43").Value2 = "Sheff Utd"

# --- Renumber column A for the rows that shifted down so the sequence
#     0..42 remains consecutive (these cells already carried their data,
#     only the running index in column A needs to be refreshed) ---
$ws.Range("A39").Value2 = 37
$ws.Range("A40").Value2 = 38
$ws.Range("A42").Value2 = 40
$ws.Range("A44").Value2 = 42
